$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "NA" value that used to be in C6 (it moves down to the new row)
$ws.Range("C6").Value = ""

# Add the new row 7 with the latest "no results" entry.
# Keep the date column formatted as text so it stays a literal string
# like the existing rows instead of becoming an Excel date serial value.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-03-07"
$ws.Range("B7").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C7").Value = "NA"
$ws.Range("D7").Value = 1
